$d = $word.ActiveDocument

# --- Requirement ID cells -------------------------------------------------
# Each ID was originally typed as two separate runs ("F_" then
# "Cust_Req_NNN"). Running Find/Replace over the full id (replacing it with
# itself) re-types the cell as a single run, which merges the two runs into
# one "F_Cust_Req_NNN" run - matching the cleaned-up document.
$ids = @(
    "F_Cust_Req_001",
    "F_Cust_Req_002",
    "F_Cust_Req_003",
    "F_Cust_Req_004",
    "F_Cust_Req_005",
    "F_Cust_Req_006",
    "F_Cust_Req_007",
    "F_Cust_Req_008",
    "F_Cust_Req_09",
    "F_Cust_Req_010",
    "F_Cust_Req_011"
)
foreach ($id in $ids) {
    $d.Content.Find.Execute($id, $true, $false, $false, $false, $false, $true, 1, $false, $id, 2) | Out-Null
}

# --- Description cells -----------------------------------------------------
# Admin features are scoped down from ADD/DELETE/UPDATE to ADD-only.
$d.Content.Find.Execute("Admin Feature - ADD/DELETE/UPDATE users.", $true, $false, $false, $false, $false, $true, 1, $false, "Admin Feature - ADD users.", 2) | Out-Null
$d.Content.Find.Execute("Admin Feature - ADD/DELETE/UPDATE  restaurants.", $true, $false, $false, $false, $false, $true, 1, $false, "Admin Feature - ADD restaurants.", 2) | Out-Null
$d.Content.Find.Execute("Admin Feature - ADD/DELETE menus ", $true, $false, $false, $false, $false, $true, 1, $false, "Admin Feature - ADD menus ", 2) | Out-Null
$d.Content.Find.Execute("Admin Feature - ADD/DELETE/UPDATE  promotions.", $true, $false, $false, $false, $false, $true, 1, $false, "Admin Feature - ADD promotions.", 2) | Out-Null

# --- Stale cursor bookmark ---------------------------------------------------
# "_GoBack" marks where the previous author last edited; it's stale now that
# the edits above moved on, so drop it if this host exposes it for removal.
try {
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
} catch {
    # Hidden-bookmark deletion isn't supported on every host; ignore.
}
